$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row at position 13 (shifts rows 13-23 down to 14-24) ---
# This new row receives "Docentes responsaveis:" value (previously
# misplaced in row 10 / row 18) under column B/C.
$ws.Rows.Item(13).Insert()

# Copy formatting (style only) from an existing B-column / C-column cell so
# the new row's cells reuse the existing style entries instead of creating
# new ones, then set their text.
$ws.Range("B3").Copy() | Out-Null
$ws.Range("B13").PasteSpecial(-4122) | Out-Null
$ws.Range("B13").Value = '5816812 - João Paulo Alves Silva'

$ws.Range("C3").Copy() | Out-Null
$ws.Range("C13").PasteSpecial(-4122) | Out-Null
$ws.Range("C13").Value = '5816812 - João Paulo Alves Silva'

# The new row has no content in column A (it's a continuation row under
# "Docentes responsaveis:" from row 12) - clear the cell Insert() created.
$ws.Range("A13").Clear()

# --- Update the cell contents that changed text (rows shifted down by the
# insert above; indices below refer to POST-insert row numbers) ---
$ws.Range("B10").Value = '1 - Estimular no aluno a capacidade de atuar como "engenheiro", no sentido de buscar soluções para o desenvolvimento de um processo químico, através da integração dos conhecimentos adquiridos em cada uma das áreas específicas da Engenharia Química. 2  Estimular o trabalho em equipe e a interação entre grupos'
$ws.Range("C10").Value = '1 - Estimular no aluno a capacidade de atuar como "engenheiro", no sentido de buscar soluções para o desenvolvimento de um processo químico, através da integração dos conhecimentos adquiridos em cada uma das áreas específicas da Engenharia Química. 2  Estimular o trabalho em equipe e a interação entre grupos'
$ws.Range("B14").Value = '1 - Visão Integradora da Engenharia Química. 2  Projetos Multidisciplinares integradores do Conhecimento em Engenharia. 3 - Desenvolvimento de projetos multidisciplinares da Indústria Química. 4 - Seminários: Apresentação e discussão dos resultados. 5. Relatório Final.'
$ws.Range("C14").Value = '1 - Visão Integradora da Engenharia Química. 2  Projetos Multidisciplinares integradores do Conhecimento em Engenharia. 3 - Desenvolvimento de projetos multidisciplinares da Indústria Química. 4 - Seminários: Apresentação e discussão dos resultados. 5. Relatório Final.'
$ws.Range("B16").Value = '1 - Visão Integradora da Engenharia Química: Análise e otimização de condições de processo de conjunto de equipamentos como reatores, trocadores de calor, sistemas de separação, entre outros.
2 - Projetos Multidisciplinares integradores do Conhecimento em Engenharia: Análise e otimização de instalações industriais.
3 - Desenvolvimento de projetos multidisciplinares da Indústria Química: Desenvolvimento de projetos visando a concepção de uma instalação industrial de uma planta química; Desenvolvimento de projetos visando o levantamento de dados e a otimização de um processo químico.   
4 - Seminários: Apresentação e discussão dos resultados. 
5 - Relatório Final'
$ws.Range("C16").Value = '1 - Visão Integradora da Engenharia Química: Análise e otimização de condições de processo de conjunto de equipamentos como reatores, trocadores de calor, sistemas de separação, entre outros.
2 - Projetos Multidisciplinares integradores do Conhecimento em Engenharia: Análise e otimização de instalações industriais.
3 - Desenvolvimento de projetos multidisciplinares da Indústria Química: Desenvolvimento de projetos visando a concepção de uma instalação industrial de uma planta química; Desenvolvimento de projetos visando o levantamento de dados e a otimização de um processo químico.   
4 - Seminários: Apresentação e discussão dos resultados. 
5 - Relatório Final'
$ws.Range("B19").Value = 'Provas escritas e Apresentação de Trabalhos'
$ws.Range("C19").Value = 'Provas escritas e Apresentação de Trabalhos'
$ws.Range("B20").Value = 'A nota será composta por ao menos uma prova escrita e trabalhos realizados e apresentados durante o semestre. O peso de cada atividade será definido segundo critérios do professor.'
$ws.Range("C20").Value = 'A nota será composta por ao menos uma prova escrita e trabalhos realizados e apresentados durante o semestre. O peso de cada atividade será definido segundo critérios do professor.'
$ws.Range("B21").Value = 'Média Final = (N + Prova Recuperação)/2'
$ws.Range("C21").Value = 'Média Final = (N + Prova Recuperação)/2'
$ws.Range("B22").Value = 'PERLINGEIRO, Carlos A. G. Engenharia de processos: análise, simulação, otimização e síntese de processos químicos.  Editora Blucher, 2005.
TURTON, BAILIE; WHITING; SHAEIWITZ  Analysis, Synthesis, and Design of Chemical Processes. 3. Ed. LTC Editora, 2005.
COULSON, J. M.; RICHARDSON, J.F. Chemical Engineering Design: Chemical Engineering Volume 6. Editora Fourth, 2005.
HIMMELBLAU, David M. Engenharia química princípios e cálculos. LTC Editora, 2006.
FELDER, R.M; Rousseau, R.W. Princípios elementares dos processos químicos. LTC Editora, 2005.
HOUGEN, O.A.; WATSON, K. M.; RAGATZ, R.A. Princípios dos processos químicos. Lopes da Silva Editora, 2005. v. 1 
CUTLIP, M.B.; SACHAM, M. Problem solving in chemical and biochemical engineering with POLYMATHTM, Excel and MATLAB®. Prentice-Hall, 2008.'
$ws.Range("C22").Value = 'PERLINGEIRO, Carlos A. G. Engenharia de processos: análise, simulação, otimização e síntese de processos químicos.  Editora Blucher, 2005.
TURTON, BAILIE; WHITING; SHAEIWITZ  Analysis, Synthesis, and Design of Chemical Processes. 3. Ed. LTC Editora, 2005.
COULSON, J. M.; RICHARDSON, J.F. Chemical Engineering Design: Chemical Engineering Volume 6. Editora Fourth, 2005.
HIMMELBLAU, David M. Engenharia química princípios e cálculos. LTC Editora, 2006.
FELDER, R.M; Rousseau, R.W. Princípios elementares dos processos químicos. LTC Editora, 2005.
HOUGEN, O.A.; WATSON, K. M.; RAGATZ, R.A. Princípios dos processos químicos. Lopes da Silva Editora, 2005. v. 1 
CUTLIP, M.B.; SACHAM, M. Problem solving in chemical and biochemical engineering with POLYMATHTM, Excel and MATLAB®. Prentice-Hall, 2008.'
